$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for "RA2Mx, RA3Mx" (2 Mohm resistor), based on the
#     existing "R1Mx, R3Mx" row (row 2), which is the closest sibling part. ---
$ws.Range("A2:J2").Copy()
$ws.Range("A3:J3").Insert()

# Fix up the inserted row's content for the new 2 Mohm resistor part.
$ws.Range("A3").Value = 16
$ws.Range("B3").Value = "Surface Mount Chip Resistor, Thick Film, AEC-Q200 CRCW Series, 2 Mohm, 100 mW, ± 1%, 75 V"
$ws.Range("C3").Value = "Common passives"
$ws.Range("D3").Value = "VISHAY DALE"
$ws.Range("E3").Value = "CRCW06032M00FKEA"
$ws.Range("F3").Value = "Newark"
$ws.Range("G3").Value = "52K8249"
$ws.Range("H3").Value = 0.004
$ws.Range("I3").Value = "RA2Mx, RA3Mx"

# Match the author's original formatting exactly: the Mfg-part/Vendor-part
# cells on this particular row end up with the worksheet default style
# (no explicit override), not the copied style from row 2.
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Style = "Normal"

# --- Insert a new row for "C1Mx" (ceramic capacitor), based on the
#     existing "R2Mx" row (now row 4), another "Common passives" part. ---
$ws.Range("A4:J4").Copy()
$ws.Range("A7:J7").Insert()

# Fix up the inserted row's content for the new ceramic capacitor part.
$ws.Range("A7").Value = 8
$ws.Range("B7").Value = "Multilayer Ceramic Capacitor, VJ.W1BC Series, 0.1 uF, 10%, X7R, 50 V, 0603 [1608 Metric]"
$ws.Range("C7").Value = "Common passives"
$ws.Range("D7").Value = "VISHAY "
$ws.Range("E7").Value = "VJ0603Y104KXACW1BC"
$ws.Range("F7").Value = "Newark"
$ws.Range("G7").Value = "52X6485"
$ws.Range("H7").Value = 0.012
$ws.Range("I7").Value = "C1Mx"

# --- Rename the two older designators that were previously "RAx " / "RBx"
#     (now shifted down to rows 5 and 6) to the new naming scheme. ---
$ws.Range("I5").Value = "RA1Mx "
$ws.Range("I6").Value = "RBMx"

# Restore the selection to match the author's final cursor position.
$ws.Range("I4").Select()
